$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.17344416360779746
$ws.Range("B1").Value = 0.17308299752210132
$ws.Range("A2").Value = -0.15097570217283796
$ws.Range("B2").Value = 0.14956444850945338
$ws.Range("A3").Value = -0.085066512183686882
$ws.Range("B3").Value = 0.084764296479036716
$ws.Range("A4").Value = -0.076764296486230066
$ws.Range("B4").Value = 0.076504243210592904
$ws.Range("A5").Value = -0.073504243214411957
$ws.Range("B5").Value = 0.072642655309546278
$ws.Range("A6").Value = -0.027544552940181433
$ws.Range("B6").Value = 0.027340411612314597
$ws.Range("A7").Value = -0.017340411621998619
$ws.Range("B7").Value = 0.017304228794714049
$ws.Range("A8").Value = -0.0073042288045357395
$ws.Range("B8").Value = 0.0072712450972289844
$ws.Range("A9").Value = -0.029652686314052712
$ws.Range("B9").Value = 0.029415650707663055
$ws.Range("A10").Value = -0.027415650711974493
$ws.Range("B10").Value = 0.027398661292025395
$ws.Range("A11").Value = -0.024398661297137636
$ws.Range("B11").Value = 0.024371401380640201
$ws.Range("A12").Value = -0.02087140138622523
$ws.Range("B12").Value = 0.020674454669425746
$ws.Range("A13").Value = -0.017174454675359385
$ws.Range("B13").Value = 0.01708440816031942
$ws.Range("A14").Value = -0.0090844081696506152
$ws.Range("B14").Value = 0.0090546583403030567
$ws.Range("A15").Value = -0.0080546583445357811
$ws.Range("B15").Value = 0.0080354172836258897
$ws.Range("A16").Value = -0.0060354172886598612
$ws.Range("B16").Value = 0.0060038638796511279
$ws.Range("A17").Value = -0.004003863884768144
$ws.Range("B17").Value = 0.0039999999933897357
$ws.Range("A18").Value = -0.016107295353847917
$ws.Range("B18").Value = 0.016092302453156293
$ws.Range("A19").Value = -0.012092302456267578
$ws.Range("B19").Value = 0.012017412489313539
$ws.Range("A20").Value = -0.0080174124926184476
$ws.Range("B20").Value = 0.0080057542425500117
$ws.Range("A21").Value = -0.0040057542459166484
$ws.Range("B21").Value = 0.0039999999966244815
$ws.Range("A22").Value = -0.060497936329708679
$ws.Range("B22").Value = 0.060177117954497206
$ws.Range("A23").Value = -0.055177117959448019
$ws.Range("B23").Value = 0.054553841609023657
$ws.Range("A24").Value = -0.020100295560072645
$ws.Range("B24").Value = 0.019999999983227212
$ws.Range("A25").Value = -0.043098102373116731
$ws.Range("B25").Value = 0.04307264166294722
$ws.Range("A26").Value = -0.040572641667205644
$ws.Range("B26").Value = 0.0405434995453966
$ws.Range("A27").Value = -0.038043499549748283
$ws.Range("B27").Value = 0.037892192193196728
$ws.Range("A28").Value = -0.021034370776651379
$ws.Range("B28").Value = 0.021013133223262592
$ws.Range("A29").Value = -0.01401313323161002
$ws.Range("B29").Value = 0.014005828092499506
$ws.Range("A30").Value = 0.045994171860085142
$ws.Range("B30").Value = -0.046211096032697796
$ws.Range("A31").Value = 0.053211096024803339
$ws.Range("B31").Value = -0.053340177298547786
$ws.Range("A32").Value = 0.059743611227789373
$ws.Range("B32").Value = -0.060024928256360965
